$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Membrillo" (Feria Lagunitas de Puerto Montt).
# It belongs right after the existing row for date 44309 (row 105), so insert a new
# row there and push everything else down by one (old row 105 -> 106, ..., old row 133 -> 134).
$ws.Rows(105).Insert()

$ws.Cells.Item(105, 1).Value = 4
$ws.Cells.Item(105, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(105, 3).Value = "Los Lagos"
$ws.Cells.Item(105, 4).Value = 45016
$ws.Cells.Item(105, 5).Value = 10
$ws.Cells.Item(105, 6).Value = "Fruta"
$ws.Cells.Item(105, 7).Value = 100104
$ws.Cells.Item(105, 8).Value = "Frutos de pepita"
$ws.Cells.Item(105, 9).Value = 100104003
$ws.Cells.Item(105, 10).Value = "Membrillo"
$ws.Cells.Item(105, 11).Value = "Champion"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 400
$ws.Cells.Item(105, 14).Value = 15000
$ws.Cells.Item(105, 15).Value = 16000
$ws.Cells.Item(105, 16).Value = 15500
$ws.Cells.Item(105, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(105, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(105, 19).Value = 861
$ws.Cells.Item(105, 20).Value = 18
